# chore: fix typo in slide
#
# Slide 14 ("Mart Management System" conclusion slide):
#   1. Fix a Khmer typo in the bullet text: "ទិន្ន៏យ" -> "ទិន្ន័យ"
#      (wrong combining sign U+17CF replaced with the correct U+17D0).
#   2. Move the diagram picture behind the subtitle placeholder in the
#      shape z-order (the subtitle shape used to be painted after/above
#      the picture; now it is painted before it).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(14)

# --- 1. Fix the typo -------------------------------------------------
$body = $s.Shapes.Item(3).TextFrame.TextRange
$para = $body.Paragraphs(6, 1)
$para.Text = "classes ផ្សេងៗ និងការទំនាក់ទំនងរវាងតារាងទិន្ន័យ (Entity Relationship Diagram)"

# --- 2. Re-order the picture / subtitle shapes ------------------------
# Shape 4 is the diagram picture, Shape 5 is the subtitle textbox.
# Bring the picture forward one step so it swaps places with the
# subtitle that currently sits in front of it.
$picture = $s.Shapes.Item(4)
$picture.ZOrder(2)

# Keep the Google Slides style shape names lined up with their new slot
# (the subtitle now occupies the picture's old slot, and vice versa).
$subtitle = $s.Shapes.Item(4)
$picture = $s.Shapes.Item(5)
$subtitle.Name = "Google Shape;192;p26"
$picture.Name = "Google Shape;193;p26"
